$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column O (year 2021) mirroring column N (year 2020) for every
# data row, copying N's formatting into O first (so number formats /
# borders / fonts match) and then overwriting the value/content that
# differs for 2021.

# Row 3: plain bordered spacer cell - just copy format, no value to change.
$ws.Range("N3").Copy($ws.Range("O3"))

# Row 4: header year value.
$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = 2021

# Row 5: percentage value.
$ws.Range("N5").Copy($ws.Range("O5"))
$ws.Range("O5").Value = 97

# Row 6.
$ws.Range("N6").Copy($ws.Range("O6"))
$ws.Range("O6").Value = 96.2

# Row 7.
$ws.Range("N7").Copy($ws.Range("O7"))
$ws.Range("O7").Value = 62.7

# Row 8.
$ws.Range("N8").Copy($ws.Range("O8"))
$ws.Range("O8").Value = 100

# Row 9.
$ws.Range("N9").Copy($ws.Range("O9"))
$ws.Range("O9").Value = 100

# Row 10: textual "-" placeholder (shared string already used by N10).
$ws.Range("N10").Copy($ws.Range("O10"))
$ws.Range("O10").Value = "-"

# Row 11.
$ws.Range("N11").Copy($ws.Range("O11"))
$ws.Range("O11").Value = 100

# Row 12.
$ws.Range("N12").Copy($ws.Range("O12"))
$ws.Range("O12").Value = 57.9

# Row 13.
$ws.Range("N13").Copy($ws.Range("O13"))
$ws.Range("O13").Value = 100

# Row 14: textual "-" placeholder.
$ws.Range("N14").Copy($ws.Range("O14"))
$ws.Range("O14").Value = "-"

# Match the saved selection state from the authored workbook.
$ws.Range("O17").Select()
